# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "226.99") are forced back
# to text (matching the original inlineStr cell type / preserving trailing zeros)
# by setting an explicit text NumberFormat before writing the value.

$ws.Range('D2').Value = '34.614.18'
$ws.Range('D3').Value = '1.794.58'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.99'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.91'
$ws.Range('E8').Value = '  +3.61%  '
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0695'
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '2.051.91'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').Value = '1.790.03'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.637'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '34.564.89'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.29'
$ws.Range('E17').Value = '  +2.65%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.82'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '248.37'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('D20').Value = '0.0₃0802'
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.30'
$ws.Range('E21').Value = '  +2.91%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '166.31'
$ws.Range('E25').Value = '  +2.82%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.30'
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.57'
$ws.Range('E27').Value = '  +1.52%  '
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.12'
$ws.Range('E30').Value = '  +13.34%  '
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0525'
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.24'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('E34').Value = '  +2.59%  '
$ws.Range('D35').Value = '1.428.11'
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('E36').Value = '  +6.52%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.673'
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0193'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.06'
$ws.Range('E39').Value = '  +1.83%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '85.28'
$ws.Range('E40').Value = '  +6.48%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.935'
$ws.Range('E42').Value = '  +1.32%  '
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.72'
$ws.Range('E44').Value = '  +1.16%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0528'
$ws.Range('E45').Value = '  +3.84%  '
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').Value = '1.952.33'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '106.12'
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('E51').Value = '  -6.14%  '
